# Updated cryptos list data (price + 1h volume change columns) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to remain a text value even when the new text looks like
    # a plain number (e.g. "173.70"), which Excel would otherwise coerce to a
    # numeric 173.7 and silently drop the trailing zero / text formatting.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '67.343.06'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '3.519.02'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '596.42'
$ws.Range('E5').Value = '  +0.96%  '
Set-TextValue 'D6' '173.70'
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('E7').Value = '  +0.00%  '
Set-TextValue 'D8' '0.595'
$ws.Range('E8').Value = '  +2.27%  '
$ws.Range('E9').Value = '  +6.29%  '
$ws.Range('E10').Value = '  -0.36%  '
Set-TextValue 'D11' '0.438'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').Value = '4.125.19'
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('E13').Value = '  +0.08%  '
Set-TextValue 'D14' '29.26'
$ws.Range('E14').Value = '  +3.63%  '
$ws.Range('D15').Value = '67.247.18'
$ws.Range('E15').Value = '  +0.82%  '
Set-TextValue 'D16' '0.0000181'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').Value = '3.552.73'
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('E19').Value = '  +1.66%  '
Set-TextValue 'D20' '396.30'
$ws.Range('E20').Value = '  +2.06%  '
Set-TextValue 'D21' '8.05'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  +0.10%  '
Set-TextValue 'D24' '0.540'
$ws.Range('E24').Value = '  +1.31%  '
Set-TextValue 'D25' '0.0000123'
$ws.Range('E25').Value = '  +0.56%  '
Set-TextValue 'D26' '10.26'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('E29').Value = '  -0.56%  '
Set-TextValue 'D30' '1.48'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  +0.57%  '
Set-TextValue 'D32' '23.90'
$ws.Range('E32').Value = '  +1.45%  '
Set-TextValue 'D33' '7.41'
$ws.Range('E33').Value = '  -0.23%  '
Set-TextValue 'D34' '1.68'
$ws.Range('E34').Value = '  +4.04%  '
Set-TextValue 'D35' '163.11'
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('E37').Value = '  +0.99%  '
Set-TextValue 'D38' '7.03'
$ws.Range('E38').Value = '  +5.84%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D39' '4.72'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D40' '0.0753'
$ws.Range('E40').Value = '  +1.07%  '
Set-TextValue 'D41' '26.66'
$ws.Range('E41').Value = '  +1.18%  '
Set-TextValue 'D42' '27.32'
$ws.Range('E42').Value = '  +2.94%  '
$ws.Range('D43').Value = '2.842.00'
$ws.Range('E43').Value = '  +0.79%  '
Set-TextValue 'D44' '2.60'
$ws.Range('E44').Value = '  +2.82%  '
Set-TextValue 'D45' '43.05'
$ws.Range('E45').Value = '  +0.14%  '
Set-TextValue 'D46' '0.0306'
$ws.Range('E46').Value = '  -1.34%  '
Set-TextValue 'D47' '340.29'
$ws.Range('E47').Value = '  -4.82%  '
$ws.Range('E48').Value = '  +0.28%  '
Set-TextValue 'D49' '34.60'
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('E50').Value = '  +0.33%  '
Set-TextValue 'D51' '0.852'
$ws.Range('E51').Value = '  -0.54%  '
